# Add data for 2022-01-06: update "through 12-28" -> "through 12-29" labels
# and refresh December / Total figures for the 2021 column (and prior years'
# December counts, which were also revised).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet and update the title label.
$ws.Name = "Through 2021-12-29"

# Update the shared-string label for December in column A (row 13).
$ws.Range("A13").Value = "December (through 12-29)"

# Update December row (row 13) values for each year column (B:H = 2015-2021).
$ws.Range("B13").Value = 44
$ws.Range("C13").Value = 93
$ws.Range("D13").Value = 111
$ws.Range("E13").Value = 72
$ws.Range("F13").Value = 62
$ws.Range("G13").Value = 136
$ws.Range("H13").Value = 182

# Update Total row (row 14) values for each year column (B:H = 2015-2021).
$ws.Range("B14").Value = 335
$ws.Range("C14").Value = 656
$ws.Range("D14").Value = 932
$ws.Range("E14").Value = 754
$ws.Range("F14").Value = 596
$ws.Range("G14").Value = 1400
$ws.Range("H14").Value = 1825
